$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed/updated) date column C for all data rows (2-70)
# from serial date 45207 (2023-10-08) to 45208 (2023-10-09).
$ws.Range("C2:C70").Value = 45208
